# Add the 4 May 2020 ("4 Mayıs 2020") row of data to the "data" sheet.
# The sheet holds a structured Table (Table3) over A1:E53; we grow it by
# one row (to A1:E54) and populate the new row with the day's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$lo = $ws.ListObjects.Item("Table3")

# Appending a ListRow extends the table (and its AutoFilter/dimension)
# by one row automatically, inheriting the existing column formatting.
$newRow = $lo.ListRows.Add()

# date, test, case, death, recovered
$ws.Range("A54").Value = 43955
$ws.Range("B54").Value = 35771
$ws.Range("C54").Value = 1614
$ws.Range("D54").Value = 64
$ws.Range("E54").Value = 5015
